# Chap05: Correction herve done
#
# 1) Refresh the cached "datetimeFigureOut" date field (Insert > Header &
#    Footer > Date and time > Update automatically) from 6/22/2017 to
#    11/7/2017 everywhere the placeholder lives: the slide master, every
#    slide layout, and (best effort) the notes master.
# 2) Nudge the "(a)" / "(b)" caption textboxes on the slide and bump their
#    font size to 24pt.

$p = $ppt.ActivePresentation

$newDate = "11/7/2017"
$ppDatePlaceholder = 16   # ppPlaceholderDate

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq $ppDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j).Shapes
}

# Notes master (best effort - not all hosts expose a writable notes master).
try {
    Update-DatePlaceholder $p.NotesMaster.Shapes
} catch {
}

# --- (a) / (b) caption textboxes on slide 1 -----------------------------

$slide = $p.Slides.Item(1)

$shapeA = $slide.Shapes.Item("ZoneTexte 17")
$shapeA.Left = 41.07811173622049
$shapeA.Top = 5.693385926771654
$shapeA.Width = 41.425669291338586
$shapeA.Height = 36.351575903149616
$shapeA.TextFrame.TextRange.Font.Size = 24

$shapeB = $slide.Shapes.Item("ZoneTexte 19")
$shapeB.Left = 324.57417322834647
$shapeB.Top = 102.0820503440943
$shapeB.Width = 42.81409448818898
$shapeB.Height = 36.351575903149616
$shapeB.TextFrame.TextRange.Font.Size = 24
